# Apply transcript renaming edits described in the commit:
#   "Renamed few trancripts. Updated the DataSheet"
#
# Changes:
#   - Column D ("Speaker") cells containing "Student" -> "S"
#     for rows: 8, 9, 11, 32, 35, 38, 39, 41, 43, 52
#   - Column F ("Teacher Tag") cells containing
#     "3 - getting students to relate" -> "3 - getting SS to relate"
#     for rows: 13, 54, 69

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$studentRows = @(8, 9, 11, 32, 35, 38, 39, 41, 43, 52)
foreach ($r in $studentRows) {
    $cell = $ws.Range("D$r")
    if ($cell.Value2 -eq "Student") {
        $cell.Value = "S"
    }
}

$tagRows = @(13, 54, 69)
foreach ($r in $tagRows) {
    $cell = $ws.Range("F$r")
    if ($cell.Value2 -eq "3 - getting students to relate") {
        $cell.Value = "3 - getting SS to relate"
    }
}
